$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D6 was a text value ("4:36") stored as a shared string; it should become
# the plain number 436. E6 keeps its original text ("17:05") - after D6's
# old string is removed from the shared-string table, Excel will simply
# re-use/compact the table automatically.
$ws.Range("D6").Value = 436
$ws.Range("E6").Value = "17:05"
